# Updates the cryptos price list: refreshed prices/volumes and a few
# re-ranked coin rows (see commit message / xml diff for details).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.418.80"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").Value = "2.427.90"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.14"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.07"
$ws.Range("E6").Value = "  +5.83%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.169"
$ws.Range("E9").Value = "  +8.82%  "
$ws.Range("D10").Value = "2.429.61"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.334"
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.67"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000178"
$ws.Range("E14").Value = "  +6.50%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "69.203.02"
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "2.873.93"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.94"
$ws.Range("E17").Value = "  +5.69%  "
$ws.Range("D18").Value = "2.423.03"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.81"
$ws.Range("E19").Value = "  +5.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.32"
$ws.Range("E20").Value = "  +4.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.14"
$ws.Range("E21").Value = "  +6.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.87"
$ws.Range("E22").Value = "  +3.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.97"
$ws.Range("E23").Value = "  +7.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.04"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.81"
$ws.Range("E26").Value = "  +6.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.49"
$ws.Range("E27").Value = "  +7.20%  "
$ws.Range("D28").Value = "2.552.46"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "0.0₃0849"
$ws.Range("E30").Value = "  +8.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.39"
$ws.Range("E31").Value = "  +6.52%  "
$ws.Range("E32").Value = "  +11.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "453.39"
$ws.Range("E33").Value = "  +9.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  +2.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.02"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.10"
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.111"
$ws.Range("E38").Value = "  +6.90%  "
$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.25"
$ws.Range("E40").Value = "  +3.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.303"
$ws.Range("E41").Value = "  +4.72%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.41"
$ws.Range("E42").Value = "  +5.48%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.52"
$ws.Range("E43").Value = "  +5.55%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "37.81"
$ws.Range("E44").Value = "  +1.94%  "
$ws.Range("E45").Value = "  +4.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.09"
$ws.Range("E46").Value = "  +9.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.56"
$ws.Range("E47").Value = "  +6.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.39"
$ws.Range("E48").Value = "  +4.73%  "
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.489"
$ws.Range("E50").Value = "  +4.38%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.561"
$ws.Range("E51").Value = "  +2.08%  "
